$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "85÷7=12, 1"
$t.Cell(1,2).Range.Text = "19÷6=3, 1"
$t.Cell(1,3).Range.Text = "50÷3=16, 2"
$t.Cell(1,4).Range.Text = "20÷8=2, 4"
$t.Cell(1,5).Range.Text = "78÷5=15, 3"
$t.Cell(5,1).Range.Text = "14÷2=7, 0"
$t.Cell(5,2).Range.Text = "45÷8=5, 5"
$t.Cell(5,3).Range.Text = "90÷9=10, 0"
$t.Cell(5,4).Range.Text = "10÷9=1, 1"
$t.Cell(5,5).Range.Text = "79÷9=8, 7"
$t.Cell(9,1).Range.Text = "96÷5=19, 1"
$t.Cell(9,2).Range.Text = "30÷2=15, 0"
$t.Cell(9,3).Range.Text = "11÷7=1, 4"
$t.Cell(9,4).Range.Text = "24÷5=4, 4"
$t.Cell(9,5).Range.Text = "25÷9=2, 7"
$t.Cell(13,1).Range.Text = "21÷5=4, 1"
$t.Cell(13,2).Range.Text = "93÷5=18, 3"
$t.Cell(13,3).Range.Text = "84÷3=28, 0"
$t.Cell(13,4).Range.Text = "92÷5=18, 2"
$t.Cell(13,5).Range.Text = "68÷3=22, 2"
$t.Cell(17,1).Range.Text = "70÷4=17, 2"
$t.Cell(17,2).Range.Text = "78÷5=15, 3"
$t.Cell(17,3).Range.Text = "26÷5=5, 1"
$t.Cell(17,4).Range.Text = "22÷9=2, 4"
$t.Cell(17,5).Range.Text = "98÷2=49, 0"
